$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1909554164245861
$ws.Range("D2").Value = 0.04649090588771543
$ws.Range("E2").Value = 0.1198196255914432
$ws.Range("F2").Value = 0.9041384715165393
$ws.Range("G2").Value = 0.002458797528885286
$ws.Range("I2").Value = 0.9159076666920143
$ws.Range("K2").Value = 1.332658977748679
$ws.Range("L2").Value = 0.1470131957561911
$ws.Range("M2").Value = 0.3985214942162756
$ws.Range("O2").Value = 3.219270872553608

$ws.Range("C3").Value = 0.1884995670253033
$ws.Range("D3").Value = 0.04442701427139184
$ws.Range("E3").Value = 0.1199813651475345
$ws.Range("F3").Value = 0.9107671145248588
$ws.Range("G3").Value = 0.002461696270925397
$ws.Range("I3").Value = 0.9243026318311678
$ws.Range("K3").Value = 1.1814544102985
$ws.Range("L3").Value = 0.1483474949222447
$ws.Range("M3").Value = 0.3677618158228384
$ws.Range("O3").Value = 3.257357974828835

$ws.Range("C4").Value = 0.187066929185093
$ws.Range("D4").Value = 0.04315422489717946
$ws.Range("E4").Value = 0.1201420431006177
$ws.Range("F4").Value = 0.9155409833865065
$ws.Range("G4").Value = 0.002463570348477426
$ws.Range("I4").Value = 0.9301062141144207
$ws.Range("K4").Value = 1.088394305530073
$ws.Range("L4").Value = 0.149252416021227
$ws.Range("M4").Value = 0.3489113508678798
$ws.Range("O4").Value = 3.283167235726197

$ws.Range("C5").Value = 0.1865021066187325
$ws.Range("D5").Value = 0.04263419304613336
$ws.Range("E5").Value = 0.1202229615258155
$ws.Range("F5").Value = 0.917663107118166
$ws.Range("G5").Value = 0.002464357817982894
$ws.Range("I5").Value = 0.9326342427150855
$ws.Range("K5").Value = 1.050418647405735
$ws.Range("L5").Value = 0.1496427076948095
$ws.Range("M5").Value = 0.3412391789052336
$ws.Range("O5").Value = 3.29429323864386

$ws.Range("C6").Value = 0.186409467021079
$ws.Range("D6").Value = 0.04254776109362268
$ws.Range("E6").Value = 0.1202373308267841
$ws.Range("F6").Value = 0.9180261508725067
$ws.Range("G6").Value = 0.002464490014267802
$ws.Range("I6").Value = 0.9330638606304227
$ws.Range("K6").Value = 1.044109693342079
$ws.Range("L6").Value = 0.1497088153368402
$ws.Range("M6").Value = 0.3399658109782351
$ws.Range("O6").Value = 3.296177429007898

$ws.Range("C7").Value = 0.1870592348202109
$ws.Range("D7").Value = 0.04314721702421309
$ws.Range("E7").Value = 0.120143071863211
$ws.Range("F7").Value = 0.9155688878454527
$ws.Range("G7").Value = 0.002463580872366623
$ws.Range("I7").Value = 0.930139648170659
$ws.Range("K7").Value = 1.087882363574977
$ws.Range("L7").Value = 0.1492575924745552
$ws.Range("M7").Value = 0.348807842000852
$ws.Range("O7").Value = 3.283314822251029

$ws.Range("C8").Value = 0.1900930553264146
$ws.Range("D8").Value = 0.04578044922185143
$ws.Range("E8").Value = 0.1198626613713287
$ws.Range("F8").Value = 0.9062777872423595
$ws.Range("G8").Value = 0.002459777498726323
$ws.Range("I8").Value = 0.9186674311191965
$ws.Range("K8").Value = 1.28057065687392
$ws.Range("L8").Value = 0.1474554782098885
$ws.Range("M8").Value = 0.387908333488781
$ws.Range("O8").Value = 3.231899729119249

$ws.Range("C9").Value = 0.1966372604061064
$ws.Range("D9").Value = 0.05089883059231681
$ws.Range("E9").Value = 0.1197994458216911
$ws.Range("F9").Value = 0.8936549491432118
$ws.Range("G9").Value = 0.002453063597374981
$ws.Range("I9").Value = 0.9013289781049636
$ws.Range("K9").Value = 1.656601484826979
$ws.Range("L9").Value = 0.1446016635936864
$ws.Range("M9").Value = 0.4648542415935069
$ws.Range("O9").Value = 3.150347142242566

$ws.Range("C10").Value = 0.2018055969671622
$ws.Range("D10").Value = 0.05463023381886245
$ws.Range("E10").Value = 0.1200495192319302
$ws.Range("F10").Value = 0.887810823439203
$ws.Range("G10").Value = 0.002448580216461203
$ws.Range("I10").Value = 0.8917477525420452
$ws.Range("K10").Value = 1.931668562851598
$ws.Range("L10").Value = 0.1429203122807543
$ws.Range("M10").Value = 0.5215347778372177
$ws.Range("O10").Value = 3.102237781795168

$ws.Range("C11").Value = 0.2042345971030528
$ws.Range("D11").Value = 0.05632116472678206
$ws.Range("E11").Value = 0.1202276487249776
$ws.Range("F11").Value = 0.8859008196099438
$ws.Range("G11").Value = 0.002446637208377715
$ws.Range("I11").Value = 0.8880773275879434
$ws.Range("K11").Value = 2.056526202202065
$ws.Range("L11").Value = 0.1422457720100567
$ws.Range("M11").Value = 0.5473494403248935
$ws.Range("O11").Value = 3.082928050803446

$ws.Range("C12").Value = 0.2051655458822097
$ws.Range("D12").Value = 0.05696051118056999
$ws.Range("E12").Value = 0.120304352871397
$ws.Range("F12").Value = 0.8852854854197361
$ws.Range("G12").Value = 0.00244591524575892
$ws.Range("I12").Value = 0.8867866012708916
$ws.Range("K12").Value = 2.103765522638582
$ws.Range("L12").Value = 0.1420033416039921
$ws.Range("M12").Value = 0.5571287782327232
$ws.Range("O12").Value = 3.075987342288443

$ws.Range("C13").Value = 0.2049645549808616
$ws.Range("D13").Value = 0.05682286031686345
$ws.Range("E13").Value = 0.1202874219425425
$ws.Range("F13").Value = 0.8854132029941155
$ws.Range("G13").Value = 0.002446070120094714
$ws.Range("I13").Value = 0.8870601673474781
$ws.Range("K13").Value = 2.09359357725117
$ws.Range("L13").Value = 0.1420549747444184
$ws.Range("M13").Value = 0.5550224577246468
$ws.Range("O13").Value = 3.077465609137221

$ws.Range("C14").Value = 0.2043109639305527
$ws.Range("D14").Value = 0.05637378384360403
$ws.Range("E14").Value = 0.1202337738639159
$ws.Range("F14").Value = 0.8858480306098073
$ws.Range("G14").Value = 0.002446577535708904
$ws.Range("I14").Value = 0.8879691496329087
$ws.Range("K14").Value = 2.060413456512606
$ws.Range("L14").Value = 0.1422255664151812
$ws.Range("M14").Value = 0.5481539172381957
$ws.Range("O14").Value = 3.082349583372263

$ws.Range("C15").Value = 0.2039120693373206
$ws.Range("D15").Value = 0.05609858380682198
$ws.Range("E15").Value = 0.1202021173625205
$ws.Range("F15").Value = 0.8861284414964388
$ws.Range("G15").Value = 0.002446890139118633
$ws.Range("I15").Value = 0.8885388510998098
$ws.Range("K15").Value = 2.04008418931852
$ws.Range("L15").Value = 0.1423317526403842
$ws.Range("M15").Value = 0.5439472295166041
$ws.Range("O15").Value = 3.085389567827292

$ws.Range("C16").Value = 0.2016484182519207
$ws.Range("D16").Value = 0.05451959345359825
$ws.Range("E16").Value = 0.1200391728292054
$ws.Range("F16").Value = 0.8879507350459832
$ws.Range("G16").Value = 0.002448709133148606
$ws.Range("I16").Value = 0.8920014901491982
$ws.Range("K16").Value = 1.923503141473759
$ws.Range("L16").Value = 0.1429662138239038
$ws.Range("M16").Value = 0.5198482990216746
$ws.Range("O16").Value = 3.103551638729869

$ws.Range("C17").Value = 0.2002796485883067
$ws.Range("D17").Value = 0.05354924175789932
$ws.Range("E17").Value = 0.1199556940721784
$ws.Range("F17").Value = 0.8892605854202387
$ws.Range("G17").Value = 0.002449849697911363
$ws.Range("I17").Value = 0.894302117672062
$ws.Range("K17").Value = 1.851913068655506
$ws.Range("L17").Value = 0.1433785785375647
$ws.Range("M17").Value = 0.5050718249159871
$ws.Range("O17").Value = 3.115353915431825

$ws.Range("C18").Value = 0.199499704315528
$ws.Range("D18").Value = 0.05299051121774312
$ws.Range("E18").Value = 0.1199137384479165
$ws.Range("F18").Value = 0.8900844125534988
$ws.Range("G18").Value = 0.002450514807883606
$ws.Range("I18").Value = 0.8956901332697242
$ws.Range("K18").Value = 1.810710919956819
$ws.Range("L18").Value = 0.1436242597110997
$ws.Range("M18").Value = 0.4965756782049624
$ws.Range("O18").Value = 3.122384620409235

$ws.Range("C19").Value = 0.1992368904331414
$ws.Range("D19").Value = 0.05280123104571999
$ws.Range("E19").Value = 0.1199005738399563
$ws.Range("F19").Value = 0.8903754343171997
$ws.Range("G19").Value = 0.002450741565411589
$ws.Range("I19").Value = 0.8961712062162235
$ws.Range("K19").Value = 1.796756288352015
$ws.Range("L19").Value = 0.1437089025073917
$ws.Range("M19").Value = 0.4936995402336208
$ws.Range("O19").Value = 3.124806685580495

$ws.Range("C20").Value = 0.2004245976782499
$ws.Range("D20").Value = 0.05365260074744072
$ws.Range("E20").Value = 0.1199639534724071
$ws.Range("F20").Value = 0.8891138575370192
$ws.Range("G20").Value = 0.002449727342721382
$ws.Range("I20").Value = 0.894050508302918
$ws.Range("K20").Value = 1.859536604322102
$ws.Range("L20").Value = 0.1433338017849763
$ws.Range("M20").Value = 0.5066445102061437
$ws.Range("O20").Value = 3.114072452733382

$ws.Range("C21").Value = 0.2045026376488437
$ws.Range("D21").Value = 0.05650571517429626
$ws.Range("E21").Value = 0.1202492805998254
$ws.Range("F21").Value = 0.8857173792250421
$ws.Range("G21").Value = 0.002446428121125085
$ws.Range("I21").Value = 0.8876994658271187
$ws.Range("K21").Value = 2.070160409765606
$ws.Range("L21").Value = 0.1421751064702477
$ws.Range("M21").Value = 0.5501712724764189
$ws.Range("O21").Value = 3.08090495131961

$ws.Range("C22").Value = 0.207232771212702
$ws.Range("D22").Value = 0.05836470565186858
$ws.Range("E22").Value = 0.1204896738842365
$ws.Range("F22").Value = 0.8841268988071675
$ws.Range("G22").Value = 0.002444352369007197
$ws.Range("I22").Value = 0.8841269197968487
$ws.Range("K22").Value = 2.207571880175578
$ws.Range("L22").Value = 0.1414936331074479
$ws.Range("M22").Value = 0.5786409530239212
$ws.Range("O22").Value = 3.061393786242974

$ws.Range("C23").Value = 0.2057697298680381
$ws.Range("D23").Value = 0.05737305978207274
$ws.Range("E23").Value = 0.120356440147269
$ws.Range("F23").Value = 0.8849180848811287
$ws.Range("G23").Value = 0.00244545289419115
$ws.Range("I23").Value = 0.8859806703373039
$ws.Range("K23").Value = 2.134255872259473
$ws.Range("L23").Value = 0.1418504071244477
$ws.Range("M23").Value = 0.5634442562603397
$ws.Range("O23").Value = 3.071608730101303

$ws.Range("C24").Value = 0.2003590444509911
$ws.Range("D24").Value = 0.05360587483173873
$ws.Range("E24").Value = 0.1199602005913114
$ws.Range("F24").Value = 0.8891799727693552
$ws.Range("G24").Value = 0.002449782630270513
$ws.Range("I24").Value = 0.8941640573398786
$ws.Range("K24").Value = 1.856090140914091
$ws.Range("L24").Value = 0.1433540185450326
$ws.Range("M24").Value = 0.5059335021256004
$ws.Range("O24").Value = 3.114651037628164

$ws.Range("C25").Value = 0.1948034368790275
$ws.Range("D25").Value = 0.0495191782941049
$ws.Range("E25").Value = 0.1197644635718689
$ws.Range("F25").Value = 0.89646871991674
$ws.Range("G25").Value = 0.002454800653064438
$ws.Range("I25").Value = 0.9054659318153853
$ws.Range("K25").Value = 1.55508022934788
$ws.Range("L25").Value = 0.1453008091123742
$ws.Range("M25").Value = 0.4440111268378359
$ws.Range("O25").Value = 3.170339998558234
